# Generate Report for Handback
# - Overview sheet: status changes from "Ready for handoff" to
#   "Handed back: in sync with en-US" for both locales.
# - zh-cn / de-de sheets: the handback finished cleanly, so the
#   "Latest Handback DateTime" is refreshed and the "Error Detail"
#   column is cleared.
# - Column widths are widened (status/date columns) and the now
#   mostly-empty Error Detail column is narrowed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2:F3").Value = "Handed back: in sync with en-US"
$ovw.Columns.Item(5).ColumnWidth = 29.15
$ovw.Columns.Item(6).ColumnWidth = 29.15

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("K2:K3").Value = "2016-09-09 12:14:09"
$zh.Range("P2:P3").Value = ""
$zh.Columns.Item(3).ColumnWidth = 29.15
$zh.Columns.Item(16).ColumnWidth = 12.8

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("K2:K3").Value = "2016-09-09 12:14:30"
$de.Range("P2:P3").Value = ""
$de.Columns.Item(3).ColumnWidth = 29.15
$de.Columns.Item(16).ColumnWidth = 12.8
